# Apply the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores numeric-looking prices as literal text (e.g. "60.573.81",
# "512.62") in the workbook. Excel auto-converts a plain numeric string typed
# into a General-formatted cell into a real number, which would corrupt values
# like "6.50" (-> 6.5) or multi-dot prices. Pre-formatting the column as Text
# keeps every write a literal string, matching the original inlineStr cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.501.96"
$ws.Range("E2").Value = "  +6.31%  "
$ws.Range("D3").Value = "2.643.12"
$ws.Range("E3").Value = "  +9.54%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "512.15"
$ws.Range("E5").Value = "  +4.67%  "
$ws.Range("D6").Value = "157.55"
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "2.683.75"
$ws.Range("E9").Value = "  +10.38%  "
$ws.Range("D10").Value = "6.49"
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("E11").Value = "  +5.06%  "
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  +3.55%  "
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "3.115.32"
$ws.Range("E14").Value = "  +9.79%  "
$ws.Range("D15").Value = "60.592.03"
$ws.Range("E15").Value = "  +6.43%  "
$ws.Range("D16").Value = "21.82"
$ws.Range("E16").Value = "  +5.14%  "
$ws.Range("E17").Value = "  +4.99%  "
$ws.Range("D18").Value = "2.681.52"
$ws.Range("E18").Value = "  +10.22%  "
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").Value = "349.06"
$ws.Range("E20").Value = "  +7.75%  "
$ws.Range("D21").Value = "10.53"
$ws.Range("E21").Value = "  +5.23%  "
$ws.Range("D22").Value = "6.21"
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "60.33"
$ws.Range("E24").Value = "  +3.30%  "
$ws.Range("E25").Value = "  +3.96%  "
$ws.Range("D26").Value = "2.786.91"
$ws.Range("E26").Value = "  +9.65%  "
$ws.Range("E27").Value = "  +3.66%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "0.0$([char]0x2083)0871"
$ws.Range("E29").Value = "  +11.35%  "
$ws.Range("D30").Value = "7.55"
$ws.Range("E30").Value = "  +3.39%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "19.62"
$ws.Range("E32").Value = "  +5.55%  "
$ws.Range("D33").Value = "157.29"
$ws.Range("E33").Value = "  +4.45%  "
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").Value = "5.74"
$ws.Range("E35").Value = "  +8.51%  "
$ws.Range("D36").Value = "4.07"
$ws.Range("E36").Value = "  +9.20%  "
$ws.Range("E37").Value = "  +4.83%  "
$ws.Range("D38").Value = "314.29"
$ws.Range("E38").Value = "  +16.85%  "
$ws.Range("E39").Value = "  +9.89%  "
$ws.Range("D40").Value = "0.856"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").Value = "3.77"
$ws.Range("E41").Value = "  +6.59%  "
$ws.Range("D42").Value = "0.836"
$ws.Range("E42").Value = "  +30.16%  "
$ws.Range("E43").Value = "  +3.84%  "
$ws.Range("D44").Value = "0.646"
$ws.Range("E44").Value = "  +8.43%  "
$ws.Range("D45").Value = "0.0576"
$ws.Range("E45").Value = "  +8.07%  "
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "20.17"
$ws.Range("E47").Value = "  +15.55%  "
$ws.Range("D48").Value = "0.990"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").Value = "  +7.47%  "
$ws.Range("D50").Value = "2.077.05"
$ws.Range("E50").Value = "  +10.53%  "
$ws.Range("D51").Value = "0.0236"
$ws.Range("E51").Value = "  +3.26%  "
